$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'25.963.99"
$ws.Range('E2').Value = "'  -2.74%  "
$ws.Range('D3').Value = "'1.836.99"
$ws.Range('E3').Value = "'  -1.22%  "
$ws.Range('D4').Value = "'1.0000"
$ws.Range('E4').Value = "'  -0.46%  "
$ws.Range('D5').Value = "'278.31"
$ws.Range('E5').Value = "'  -4.37%  "
$ws.Range('D6').Value = "'0.9998"
$ws.Range('E6').Value = "'  -0.31%  "
$ws.Range('D7').Value = "'0.5082"
$ws.Range('E7').Value = "'  -3.24%  "
$ws.Range('D8').Value = "'0.3487"
$ws.Range('E8').Value = "'  -5.72%  "
$ws.Range('D9').Value = "'44.83"
$ws.Range('E9').Value = "'  -1.33%  "
$ws.Range('D10').Value = "'0.06797"
$ws.Range('E10').Value = "'  -4.16%  "
$ws.Range('D11').Value = "'19.83"
$ws.Range('E11').Value = "'  -5.93%  "
$ws.Range('D12').Value = "'0.8015"
$ws.Range('E12').Value = "'  -8.68%  "
$ws.Range('D13').Value = "'0.07794"
$ws.Range('E13').Value = "'  -3.61%  "
$ws.Range('D14').Value = "'1.837.89"
$ws.Range('E14').Value = "'  -1.69%  "
$ws.Range('D15').Value = "'5.064"
$ws.Range('E15').Value = "'  -3.15%  "
$ws.Range('D16').Value = "'88.11"
$ws.Range('E16').Value = "'  -3.55%  "
$ws.Range('D17').Value = "'1.000"
$ws.Range('E17').Value = "'  -0.71%  "
$ws.Range('D18').Value = "'14.14"
$ws.Range('E18').Value = "'  -3.02%  "
$ws.Range('D19').Value = "'0.000008048"
$ws.Range('E19').Value = "'  -4.44%  "
$ws.Range('D20').Value = "'0.9986"
$ws.Range('E20').Value = "'  -0.38%  "
$ws.Range('D21').Value = "'26.008.90"
$ws.Range('E21').Value = "'  -2.72%  "
$ws.Range('D22').Value = "'4.769"
$ws.Range('E22').Value = "'  -2.98%  "
$ws.Range('D23').Value = "'10.03"
$ws.Range('E23').Value = "'  -4.75%  "
$ws.Range('D24').Value = "'6.192"
$ws.Range('E24').Value = "'  -1.77%  "
$ws.Range('D25').Value = "'2.357"
$ws.Range('E25').Value = "'  +5.46%  "
$ws.Range('D26').Value = "'143.21"
$ws.Range('E26').Value = "'  -1.28%  "
$ws.Range('D27').Value = "'1.664"
$ws.Range('E27').Value = "'  -3.98%  "
$ws.Range('D28').Value = "'17.15"
$ws.Range('E28').Value = "'  -3.64%  "
$ws.Range('D29').Value = "'109.54"
$ws.Range('E29').Value = "'  -2.94%  "
$ws.Range('D30').Value = "'4.353"
$ws.Range('E30').Value = "'  -6.38%  "
$ws.Range('D31').Value = "'4.268"
$ws.Range('E31').Value = "'  -6.47%  "
$ws.Range('D32').Value = "'0.08786"
$ws.Range('E32').Value = "'  -2.91%  "
$ws.Range('D33').Value = "'0.04846"
$ws.Range('E33').Value = "'  -2.24%  "
$ws.Range('D34').Value = "'1.160"
$ws.Range('E34').Value = "'  +0.96%  "
$ws.Range('D35').Value = "'0.7262"
$ws.Range('E35').Value = "'  -8.11%  "
$ws.Range('E36').Value = "'  -2.79%  "
$ws.Range('D37').Value = "'3.194"
$ws.Range('E37').Value = "'  +0.54%  "
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').Value = "'2.348"
$ws.Range('E38').Value = "'  -8.78%  "
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = "'0.01845"
$ws.Range('E39').Value = "'  -3.87%  "
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').Value = "'0.5116"
$ws.Range('E40').Value = "'  -12.77%  "
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').Value = "'0.9473"
$ws.Range('E41').Value = "'  -10.06%  "
$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D42').Value = "'116.91"
$ws.Range('E42').Value = "'  +2.27%  "
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').Value = "'6.221"
$ws.Range('E43').Value = "'  -2.89%  "
$ws.Range('B44').Value = 'Aptos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D44').Value = "'7.956"
$ws.Range('E44').Value = "'  -6.47%  "
$ws.Range('B45').Value = 'PaxDollar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D45').Value = "'0.9995"
$ws.Range('E45').Value = "'  -0.32%  "
$ws.Range('B46').Value = 'Algorand'
$ws.Range('C46').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D46').Value = "'0.1359"
$ws.Range('E46').Value = "'  -7.55%  "
$ws.Range('D47').Value = "'0.4485"
$ws.Range('E47').Value = "'  -12.90%  "
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = "'9.296"
$ws.Range('E48').Value = "'  -6.38%  "
$ws.Range('B49').Value = 'Elrond'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D49').Value = "'36.02"
$ws.Range('E49').Value = "'  -2.01%  "
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = "'0.05915"
$ws.Range('E50').Value = "'  -2.04%  "
$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').Value = "'1.483"
$ws.Range('E51').Value = "'  -7.76%  "
